$wb = $excel.ActiveWorkbook

# --- sheet "isa_template": bump template version ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.4"

# --- sheet "4COm05_Imaging": update building-block header labels ---
$wsImaging = $wb.Worksheets.Item("4COm05_Imaging")
$wsImaging.Range("B1").Value = "Characteristic [use of image]"
$wsImaging.Range("H1").Value = "Component [image processing software]"
$wsImaging.Range("Q1").Value = "Component [feature extraction software]"
$wsImaging.Range("W1").Value = "Component [image data analysis software]"
$wsImaging.Range("Z1").Value = "Output [Data]"

# tidy stray trailing whitespace in the software-name value
$wsImaging.Range("Q2").Value = "R package EBImage"
